# edit.ps1 - apply the "updated documentation and slides" commit to
# RunningWorkflows.pptx via PowerPoint COM interop.
#
# Two kinds of edits are described by the diff:
#   1. Every cached "datetimeFigureOut" field (11 slide layouts, the
#      slide master and the notes master) is re-stamped from 8/29/12 to
#      8/31/12.
#   2. The speaker notes on slide 1 are rewritten: the three runs that
#      made up "Run a canned workflow ... reads quality)," are merged,
#      the trailing "... and run the workflow." sentence is extended
#      with a new "Also, myExperiment." mention, and the final period is
#      now its own run.

$p = $ppt.ActivePresentation

$oldDate = "8/29/12"
$newDate = "8/31/12"

function Set-DatePlaceholderText {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 1a. Slide master date placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes

# 1b. Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes
}

# 1c. Notes master date placeholder (best-effort - some hosts treat the
# notes master as read-only, so guard it defensively).
try {
    $notesMaster = $p.NotesMaster
    Set-DatePlaceholderText $notesMaster.Shapes
} catch {
    Write-Output "notes master date placeholder update skipped: $_"
}

# 2. Rewrite the speaker notes paragraph on slide 1.
$slide = $p.Slides.Item(1)
$notesPage = $slide.NotesPage
$notesShapes = $notesPage.Shapes
for ($i = 1; $i -le $notesShapes.Count; $i++) {
    $shp = $notesShapes.Item($i)
    if ($shp.Name -like "Notes Placeholder*") {
        $shp.TextFrame.TextRange.Text = "Run a canned workflow in galaxy. Example: prepare an example workflow on galaxy main (e.g. basic reads quality), make students register for main (https://main.g2.bx.psu.edu) and run the workflow. Also, myExperiment."
    }
}
